$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update WVOS row (row 5): membership count and region
$ws.Range("B5").Value = 87
$ws.Range("K5").Value = "Virginia"

# Update OSNJ row (row 7): membership count and region
$ws.Range("B7").Value = 649
$ws.Range("K7").Value = "New Jersey"

# Update DSCO row (row 6): membership count and region
$ws.Range("B6").Value = 104
$ws.Range("K6").Value = "Delaware"

# Update ESHOS row (row 8): membership count (region stays "New York")
$ws.Range("B8").Value = 200

# Update the active cell selection shown in the sheet view
$ws.Range("H11").Select()
